$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, shifting existing rows 72..138 down to 73..139
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data record
$ws.Cells.Item(72, 1).Value = 4
$ws.Cells.Item(72, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(72, 3).Value = "Los Lagos"
$ws.Cells.Item(72, 4).Value = 44810
$ws.Cells.Item(72, 5).Value = 10
$ws.Cells.Item(72, 6).Value = 100112052
$ws.Cells.Item(72, 7).Value = "Albahaca"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 80
$ws.Cells.Item(72, 11).Value = 5500
$ws.Cells.Item(72, 12).Value = 5500
$ws.Cells.Item(72, 13).Value = 5500
$ws.Cells.Item(72, 14).Value = "$/docena de matas"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 917
$ws.Cells.Item(72, 17).Value = 6
$ws.Cells.Item(72, 18).Value = "Hortaliza"
